# feat: change PFAS selection criteria
# Removes the "6:2 FTSA" and "PFPeA" rows from the PFAS table and
# recalculates the sumPFAS totals row to reflect the remaining substances.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 is "6:2 FTSA" -> delete entire row (everything below shifts up).
$ws.Rows.Item(2).Delete()

# After the shift above, "PFPeA" (originally row 11) is now row 10.
$ws.Rows.Item(10).Delete()

# "sumPFAS" is now row 13 - update its totals for the remaining rows.
$ws.Range("B13").Value = 14.21
$ws.Range("C13").Value = 17.03
$ws.Range("D13").Value = 25.71
$ws.Range("E13").Value = 10.58
$ws.Range("F13").Value = 15.1
$ws.Range("G13").Value = 18.7
$ws.Range("H13").Value = 0.76
$ws.Range("I13").Value = 1.13
$ws.Range("J13").Value = 2.43
